$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new task in column B, row 10, matching the "Bad" style used by
# the rest of column B (same style class as B3:B6).
$ws.Range("B10").Value = "Add Centrifugal Forces"
$ws.Range("B10").Style = "Bad"

# Move the current selection to B11 (next empty row in column B).
$ws.Range("B11").Select()
